# feature: gerenciador empresa com verificador, bug ao escrever na planilha
#
# 1) "Empresa" sheet: fill in the first (previously blank "insert") row of
#    Tabela5 with the company's registration data entered through the
#    new company manager screen.
$wb = $excel.ActiveWorkbook

$wsEmpresa = $wb.Worksheets.Item("Empresa")
$wsEmpresa.Range("A2:D2").NumberFormat = "@"
$wsEmpresa.Range("A2").Value = "Sonda"
$wsEmpresa.Range("B2").Value = "123123"
$wsEmpresa.Range("C2").Value = "Rua 2"
$wsEmpresa.Range("D2").Value = "119999999"
$wsEmpresa.Range("A2:D2").ClearFormats()

# 2) "Controle de Acesso" sheet: append the new login-audit rows recorded by
#    the access-control / login verifier (Tabela11 grows A1:C4 -> A1:C8).
$wsAcesso = $wb.Worksheets.Item("Controle de Acesso")
$lo = $wsAcesso.ListObjects.Item("Tabela11")

$dateFormat = $wsAcesso.Range("B2").NumberFormat

$loginLog = @(
    @("0001", 45444.091481516203),
    @("0001", 45444.122418749997),
    @("0002", 45444.124294965281),
    @("0003", 45444.131129421294)
)

foreach ($entry in $loginLog) {
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range.Row
    $wsAcesso.Cells.Item($r, 1).NumberFormat = "@"
    $wsAcesso.Cells.Item($r, 1).Value = $entry[0]
    $wsAcesso.Cells.Item($r, 2).NumberFormat = $dateFormat
    $wsAcesso.Cells.Item($r, 2).Value = $entry[1]
}

$wsAcesso.Range("D5").Select()
